# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on Sheet1 to match the refreshed cryptos snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.842.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "
$ws.Range("D3").Value = "'2.409.69"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'551.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "'136.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.92%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "'0.592"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.21%  "
$ws.Range("E9").Value = "  -1.35%  "
$ws.Range("D10").Value = "'5.68"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("E11").Value = "  -0.99%  "
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "'2.839.21"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.50%  "
$ws.Range("D15").Value = "'59.780.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'2.441.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.53%  "
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").Value = "'328.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.55%  "
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'66.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("D27").Value = "'1.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("D28").Value = "'0.0₃0771"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.01%  "
$ws.Range("E29").Value = "  -2.03%  "
$ws.Range("D30").Value = "'169.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.08%  "
$ws.Range("E31").Value = "  -4.27%  "
$ws.Range("E32").Value = "  -0.39%  "
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("E35").Value = "  -0.56%  "
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("E38").Value = "  -2.05%  "
$ws.Range("D39").Value = "'320.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.24%  "
$ws.Range("E40").Value = "  -4.24%  "
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").Value = "'139.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.00%  "
$ws.Range("D43").Value = "'0.0968"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.44%  "
$ws.Range("E44").Value = "  +1.25%  "
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").Value = "'0.578"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.41%  "
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  -6.85%  "
$ws.Range("D49").Value = "'11.05"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("E50").Value = "  -3.19%  "
$ws.Range("E51").Value = "  -1.22%  "
